# Insert a new data row at spreadsheet row 887. This shifts the existing
# rows 887-967 down to 888-968, preserving all of their data, and a brand
# new row of data is placed into row 887.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(887).Insert()

$ws.Cells.Item(887, 1).Value = 8
$ws.Cells.Item(887, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(887, 3).Value = "Coquimbo"
$ws.Cells.Item(887, 4).Value = 45166
$ws.Cells.Item(887, 5).Value = 4
$ws.Cells.Item(887, 6).Value = 100112043
$ws.Cells.Item(887, 7).Value = "Pepino ensalada"
$ws.Cells.Item(887, 8).Value = "Sin especificar"
$ws.Cells.Item(887, 9).Value = "Primera"
$ws.Cells.Item(887, 10).Value = 600
$ws.Cells.Item(887, 11).Value = 7500
$ws.Cells.Item(887, 12).Value = 8500
$ws.Cells.Item(887, 13).Value = 8000
$ws.Cells.Item(887, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(887, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(887, 16).Value = 133
$ws.Cells.Item(887, 17).Value = 60
$ws.Cells.Item(887, 18).Value = "Hortaliza"
